$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row for TV7377 age 127 right after the existing TV7377/126 row (row 73 -> new row 74)
$ws.Rows.Item(74).Insert()
$ws.Range("A74").Value = "TV7377"
$ws.Range("B74").Value = 127
$ws.Range("C74").Value = 961.5384615384614
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 961.5384615384614

# Insert new row for GRF95 age 127 right after the (shifted) GRF95/126 row (now at row 146 -> new row 147)
$ws.Rows.Item(147).Insert()
$ws.Range("A147").Value = "GRF95"
$ws.Range("B147").Value = 127
$ws.Range("C147").Value = 961.5384615384614
$ws.Range("D147").Value = 0
$ws.Range("E147").Value = 961.5384615384614

# Insert new row for GRM95 age 127 at the end (after the shifted GRM95/126 row, now at row 219 -> new row 220)
$ws.Rows.Item(220).Insert()
$ws.Range("A220").Value = "GRM95"
$ws.Range("B220").Value = 127
$ws.Range("C220").Value = 961.5384615384614
$ws.Range("D220").Value = 0
$ws.Range("E220").Value = 961.5384615384614
